# One-click update from Van Paper 07:10 AM on 2025-11-10
#
# Changes applied:
#  1. A23 (Customer Name): "SLAPSTIX" -> "WOODBURY ICE"
#  2. D17 (Last Invoice Date): blank -> 11/07/2025
#  3. C23 (Prospect): "040" -> "023"
#  4. D23 (Last Invoice Date): blank -> 11/07/2025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date format used by the other "Last Invoice Date" cells in column D.
$dateFormat = "[$-409]MM\/dd\/yyyy"
$invoiceDate = Get-Date -Year 2025 -Month 11 -Day 7 -Hour 0 -Minute 0 -Second 0

# 1. Rename customer "SLAPSTIX" to "WOODBURY ICE"
$ws.Range("A23").Value = "WOODBURY ICE"

# 2. Set the Last Invoice Date for row 17 (VINCENT MANUFACTURING)
#    (vertical "top" alignment is already inherited from the blank cell's
#    existing style, so only the number format and horizontal alignment
#    need to be set to match the other date cells in column D)
$ws.Range("D17").Value = $invoiceDate
$ws.Range("D17").NumberFormat = $dateFormat
$ws.Range("D17").HorizontalAlignment = -4131

# 3. Update the Prospect code for row 23 (WOODBURY ICE) from 040 to 023
$ws.Range("C23").Value = "023"

# 4. Set the Last Invoice Date for row 23 (WOODBURY ICE)
$ws.Range("D23").Value = $invoiceDate
$ws.Range("D23").NumberFormat = $dateFormat
$ws.Range("D23").HorizontalAlignment = -4131
